$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. DSA sheet - just reselect the title cell (matches final saved selection)
# ---------------------------------------------------------------------------
$dsa = $wb.Worksheets.Item("DSA")
$dsa.Select()
$dsa.Range("A1").Select()

# ---------------------------------------------------------------------------
# 2. OOP sheet - scroll so column B is left-most, select the whole data block
# ---------------------------------------------------------------------------
$oop = $wb.Worksheets.Item("OOP")
$oop.Select()
$excel.ActiveWindow.ScrollColumn = 2
$oop.Range("A1:K10").Select()

# ---------------------------------------------------------------------------
# 3. Final sheet - select A7 (leftover selection once user tabs away)
# ---------------------------------------------------------------------------
$final = $wb.Worksheets.Item("Final")
$final.Select()
$final.Range("A7").Select()

# ---------------------------------------------------------------------------
# 4. FSAD sheet - populate the assessment tracker, exactly like the other
#    module sheets (OOP, MLFCS, DSA, ...)
# ---------------------------------------------------------------------------
$fsad = $wb.Worksheets.Item("FSAD")
$fsad.Select()

$fsad.Range("A1:G1").Merge()
$fsad.Range("A1").Value = "Object Oriented Programming"

$fsad.Range("A2").Value = "Assessment"
$fsad.Range("B2").Value = "Type"
$fsad.Range("C2").Value = "Date"
$fsad.Range("D2").Value = "Weight"
$fsad.Range("E2").Value = "Marks"
$fsad.Range("F2").Value = "Out Of"
$fsad.Range("G2").Value = "Percent"
$fsad.Range("J2").Value = "Weight"
$fsad.Range("K2").Value = "Score"

$fsad.Range("A3").Value = "Assignment 1"
$fsad.Range("B3").Value = "AS"
$fsad.Range("C3").Value = [DateTime]::new(2022,3,14)
$fsad.Range("D3").Value = 0.25
$fsad.Range("F3").Value = 100
$fsad.Range("G3").Formula = "=IFERROR(OOP_9[[#This Row],[Marks]]/OOP_9[[#This Row],[Out Of]],0)"
$fsad.Range("I3").Value = "Assignments"
$fsad.Range("J3").Formula = '=SUMIF(OOP_9[Type],"AS",OOP_9[Weight])'
$fsad.Range("K3").Formula = '=(SUMIF(OOP_9[Type],"AS",OOP_9[Percent]))/COUNTIF(OOP_9[Type],"AS")'

$fsad.Range("A4").Value = "Assignment 2"
$fsad.Range("B4").Value = "AS"
$fsad.Range("C4").Value = [DateTime]::new(2022,5,3)
$fsad.Range("D4").Value = 0.5
$fsad.Range("F4").Value = 100
$fsad.Range("G4").Formula = "=IFERROR(OOP_9[[#This Row],[Marks]]/OOP_9[[#This Row],[Out Of]],0)"
$fsad.Range("I4").Value = "Total"
$fsad.Range("J4").Formula = "=SUM(J3:J3)"
$fsad.Range("K4").Formula = "=(J3*K3)"

$fsad.Range("A5").Value = "Assignment 3"
$fsad.Range("B5").Value = "AS"
$fsad.Range("C5").Value = [DateTime]::new(2022,5,16)
$fsad.Range("D5").Value = 0.25
$fsad.Range("F5").Value = 100
$fsad.Range("G5").Formula = "=IFERROR(OOP_9[[#This Row],[Marks]]/OOP_9[[#This Row],[Out Of]],0)"

# formatting for rows 3-10 mirrors the other module sheets
$fsad.Range("A3:B10").HorizontalAlignment = -4108
$fsad.Range("C3:C10").NumberFormat = "m/d/yyyy"
$fsad.Range("C3:C10").HorizontalAlignment = -4108
$fsad.Range("D3:D10").NumberFormat = "0%"
$fsad.Range("D3:D10").HorizontalAlignment = -4108
$fsad.Range("E3:F10").NumberFormat = "0.00"
$fsad.Range("E3:F10").HorizontalAlignment = -4108
$fsad.Range("G3:G10").NumberFormat = "0.00%"
$fsad.Range("G3:G10").HorizontalAlignment = -4108
$fsad.Range("J3:K4").NumberFormat = "0.00%"

# Create the structured table over the header + 3 assessment rows
$tbl = $fsad.ListObjects.Add(1, $fsad.Range("A2:G5"), [System.Type]::Missing, 1)
$tbl.Name = "OOP_9"
$tbl.TableStyle = "TableStyleLight1"

$fsad.Range("A1").Select()
$excel.ActiveWindow.Zoom = 205

# Move the "Menu" button shape down below the new table, like on OOP/MLFCS
$shp = $fsad.Shapes.Item(1)
$shp.Left = 86.46
$shp.Top = 99.32
$shp.Width = 131.47
$shp.Height = 39.83

$fsad.Range("E9").Select()
